# Correct "IG" -> "WG" in the title of the "Tuesday" agenda/header slide.
# (The WoT IG/WG distinction: prior days' headers said "IG" but should say "WG".)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

# Sanity check: make sure we are editing the expected run before touching it.
if ($tr.Text.StartsWith("IG: Tuesday June 22")) {
    # Replace the leading "I" with "W" (keeps its own run),
    # then re-touch the following "G" (becomes its own run too),
    # leaving the rest of the text (": Tuesday June 22 (1h55m)") untouched.
    $tr.Characters(1, 1).Text = "W"
    $tr.Characters(2, 1).Text = "G"
}
